{"js": "// Replace each \"before\" arithmetic expression with its corresponding\n// \"after\" expression, preserving all run/paragraph formatting by using\n// Range.insertText(..., Word.InsertLocation.replace) on search hits.\nconst replacements = [\n  [\"821\u00d74=3284\", \"580\u00d77=4060\"],\n  [\"741\u00d79=6669\", \"691\u00d78=5528\"],\n  [\"177\u00d78=1416\", \"767\u00d74=3068\"],\n  [\"413\u00d77=2891\", \"938\u00d73=2814\"],\n  [\"645\u00d79=5805\", \"749\u00d78=5992\"],\n  [\"864\u00d78=6912\", \"689\u00d74=2756\"],\n  [\"796\u00d77=5572\", \"264\u00d75=1320\"],\n  [\"602\u00d76=3612\", \"676\u00d77=4732\"],\n  [\"547\u00d73=1641\", \"555\u00d77=3885\"],\n  [\"265\u00d78=2120\", \"584\u00d77=4088\"],\n  [\"475\u00d79=4275\", \"929\u00d74=3716\"],\n  [\"751\u00d79=6759\", \"854\u00d79=7686\"],\n  [\"454\u00d79=4086\", \"893\u00d74=3572\"],\n  [\"735\u00d75=3675\", \"913\u00d74=3652\"],\n  [\"785\u00d78=6280\", \"490\u00d78=3920\"],\n  [\"966\u00d78=7728\", \"147\u00d79=1323\"],\n  [\"915\u00d77=6405\", \"430\u00d72=860\"],\n  [\"169\u00d78=1352\", \"689\u00d74=2756\"],\n  [\"344\u00d77=2408\", \"328\u00d77=2296\"],\n  [\"200\u00d74=800\", \"508\u00d73=1524\"],\n  [\"925\u00d73=2775\", \"941\u00d78=7528\"],\n  [\"692\u00d75=3460\", \"286\u00d75=1430\"],\n  [\"252\u00d74=1008\", \"588\u00d72=1176\"],\n  [\"741\u00d73=2223\", \"695\u00d77=4865\"],\n  [\"262\u00d77=1834\", \"619\u00d73=1857\"],\n];\n\nconst body = context.document.body;\n\nfor (const [before, after] of replacements) {\n  const results = body.search(before, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`No match found for \"${before}\"`);\n  }\n\n  // Only the first occurrence should exist (each \"before\" string is\n  // unique in the document), but replace all matches defensively.\n  for (const item of results.items) {\n    item.insertText(after, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each \"before\" arithmetic expression with its corresponding\n# \"after\" expression throughout the document body, using Find/Replace\n# so existing run/paragraph formatting on each cell is left untouched.\n\n$wdReplaceAll    = 2\n$wdFindContinue  = 1\n\n$d = $word.ActiveDocument\n\n$pairs = [ordered]@{\n    \"821\u00d74=3284\" = \"580\u00d77=4060\"\n    \"741\u00d79=6669\" = \"691\u00d78=5528\"\n    \"177\u00d78=1416\" = \"767\u00d74=3068\"\n    \"413\u00d77=2891\" = \"938\u00d73=2814\"\n    \"645\u00d79=5805\" = \"749\u00d78=5992\"\n    \"864\u00d78=6912\" = \"689\u00d74=2756\"\n    \"796\u00d77=5572\" = \"264\u00d75=1320\"\n    \"602\u00d76=3612\" = \"676\u00d77=4732\"\n    \"547\u00d73=1641\" = \"555\u00d77=3885\"\n    \"265\u00d78=2120\" = \"584\u00d77=4088\"\n    \"475\u00d79=4275\" = \"929\u00d74=3716\"\n    \"751\u00d79=6759\" = \"854\u00d79=7686\"\n    \"454\u00d79=4086\" = \"893\u00d74=3572\"\n    \"735\u00d75=3675\" = \"913\u00d74=3652\"\n    \"785\u00d78=6280\" = \"490\u00d78=3920\"\n    \"966\u00d78=7728\" = \"147\u00d79=1323\"\n    \"915\u00d77=6405\" = \"430\u00d72=860\"\n    \"169\u00d78=1352\" = \"689\u00d74=2756\"\n    \"344\u00d77=2408\" = \"328\u00d77=2296\"\n    \"200\u00d74=800\"  = \"508\u00d73=1524\"\n    \"925\u00d73=2775\" = \"941\u00d78=7528\"\n    \"692\u00d75=3460\" = \"286\u00d75=1430\"\n    \"252\u00d74=1008\" = \"588\u00d72=1176\"\n    \"741\u00d73=2223\" = \"695\u00d77=4865\"\n    \"262\u00d77=1834\" = \"619\u00d73=1857\"\n}\n\nforeach ($key in $pairs.Keys) {\n    $rng = $d.Content\n    $rng.Find.Execute($key, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $pairs[$key], $wdReplaceAll)\n}\n"}
